$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 457
$ws.Range("I2").Value = 1300
$ws.Range("J2").Value = 5407
$ws.Range("K2").Value = 25
$ws.Range("L2").Value = 1469
$ws.Range("M2").Value = 97
$ws.Range("N2").Value = 956
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 19
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 61
$ws.Range("S2").Value = 494
$ws.Range("T2").Value = 881
$ws.Range("U2").Value = 67
$ws.Range("V2").Value = 8139
$ws.Range("X2").Value = 8030
$ws.Range("Y2").Value = 18
$ws.Range("Z2").Value = 129
$ws.Range("AA2").Value = 50
